$d = $word.ActiveDocument

# --- Q1 answer paragraph: replace the multi-run text with a single plain run,
# then merge with the following paragraph (the horizontal-rule "pict") so both
# live in the same <w:p>. ---
$p = $d.Paragraphs.Item(2)
$full = $d.Range($p.Range.Start, $p.Range.End)
$full.Text = "The most challenging part was creating a reliable way to store and retrieve user data, especially when ensuring that the data (income and expenses) were consistently saved and loaded correctly. I also had to deal with ensuring that the input from the user was valid and accounted for exceptions, such as entering non-numeric values."
$p = $d.Paragraphs.Item(2)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# --- Q2 answer paragraph: same pattern. ---
$p = $d.Paragraphs.Item(4)
$full = $d.Range($p.Range.Start, $p.Range.End)
$full.Text = "I tested the program by running various scenarios, such as entering invalid input (non-numeric values), and ensuring that the budget calculations and chart generation worked as expected. I am confident that the program works reliably for typical use cases, but further testing with edge cases (e.g., very large numbers, negative income) would be beneficial."
$p = $d.Paragraphs.Item(4)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# --- Q3 answer paragraph: replace text only, no pict paragraph follows. ---
$p = $d.Paragraphs.Item(6)
$full = $d.Range($p.Range.Start, $p.Range.End)
$full.Text = "One potential addition would be to add more detailed categories for expenses, such as differentiating between fixed and variable expenses. I would also like to explore generating more advanced reports or graphs, such as pie charts or line charts to track expenses over time."

# --- Remove the trailing empty paragraph at the end of the document. ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$prevEnd = $last.Range.Start
$mark = $d.Range($prevEnd - 1, $prevEnd)
$mark.Delete()
